$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: "qwswd" -> number 12345 (B1/C1 keep their original text: "eerer"/"wwww") ---
$ws.Range("A1").Value = 12345

# --- Row 8: C8 = "harika", H8 = hyperlink to abc@gmail.com ---
$ws.Range("C8").Value = "harika"
$ws.Hyperlinks.Add($ws.Range("H8"), "mailto:abc@gmail.com", "", "mailto:abc@gmail.com", "abc@gmail.com")

# --- Row 9: H9 = hyperlink "bava@123" (no scheme -> plain external target, no tooltip) ---
$ws.Hyperlinks.Add($ws.Range("H9"), "bava@123", "", "", "bava@123")
# Give it the "Followed Hyperlink" look (distinct style/font from H8)
$ws.Range("H9").Style = "Followed Hyperlink"

# --- Row 10: H10 carries the same "Hyperlink" style as H8 but holds no value ---
$ws.Range("H10").Style = "Hyperlink"

# --- Row 12: F12 = "anam" ---
$ws.Range("F12").Value = "anam"

# --- Column H formatting: outline level + width ---
$ws.Columns("H").OutlineLevel = 7
$ws.Columns("H").ColumnWidth = 9.83

# --- Final selection lands on H10 ---
$ws.Range("H10").Select() | Out-Null
